$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 185.3901149227057
$ws.Range("AD2").Value = 187.9884851080575
$ws.Range("F3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("Q3").Value = 28.28668428403349
$ws.Range("S3").Value = 150.1310015817778
$ws.Range("W3").Value = 196.8604391758524
$ws.Range("AA3").Value = 0
$ws.Range("G4").Value = 1.111717471350082
$ws.Range("S4").Value = 0
$ws.Range("AB4").Value = 167.775720771107
$ws.Range("C5").Value = 49.40431712919045
$ws.Range("I5").Value = 0
$ws.Range("L5").Value = 46.28321089717047
$ws.Range("O5").Value = 0
$ws.Range("R5").Value = 167.9153402463508
$ws.Range("Z5").Value = 100.7553394673655
$ws.Range("C6").Value = 0
$ws.Range("J6").Value = 181.5221913118031
$ws.Range("P6").Value = 0
$ws.Range("T6").Value = 193.0410977412918
$ws.Range("U6").Value = 152.856166140868
$ws.Range("Y6").Value = 159.6626230917951
$ws.Range("AB6").Value = 46.56762015541631
$ws.Range("H7").Value = 7.275451017233925
$ws.Range("L7").Value = 80.65947074088903
$ws.Range("N7").Value = 43.56036642424102
$ws.Range("R7").Value = 76.94643731847638
$ws.Range("S7").Value = 119.2290395715406
$ws.Range("U7").Value = 132.5830079260408
$ws.Range("I8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("Q8").Value = 69.9595483632778
$ws.Range("X8").Value = 60.69780122155435
$ws.Range("Z8").Value = 150.0011984146961
$ws.Range("AA8").Value = 40.15186014115548
$ws.Range("AA9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("E10").Value = 61.61995085843572
$ws.Range("P10").Value = 66.38174720444958
$ws.Range("U10").Value = 146.9720258416792
$ws.Range("Y10").Value = 122.7429672168996
$ws.Range("AA10").Value = 192.0267570975421
$ws.Range("C11").Value = 131.7461276268717
$ws.Range("I11").Value = 98.06645730313166
$ws.Range("L11").Value = 0
$ws.Range("S11").Value = 59.8127427345013
$ws.Range("V11").Value = 140.5124702930762
$ws.Range("AA11").Value = 135.2030467813716
$ws.Range("AD11").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 186.9376378245593
$ws.Range("M12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("S12").Value = 122.2396488798556
$ws.Range("U12").Value = 50.8611802002809
$ws.Range("Y12").Value = 114.9802407312519
$ws.Range("Z12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("X13").Value = 55.67383796497867
$ws.Range("AC13").Value = 133.8674284005627
$ws.Range("AD13").Value = 0
$ws.Range("AE13").Value = 139.5846115125957
$ws.Range("E14").Value = 119.9388001706701
$ws.Range("F14").Value = 176.1177051987486
$ws.Range("J14").Value = 192.5400163637157
$ws.Range("M14").Value = 64.79984791280602
$ws.Range("R14").Value = 0
$ws.Range("Y14").Value = 0
$ws.Range("AD14").Value = 147.2582676252367
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 59.13755813526858
$ws.Range("L15").Value = 63.63715947293589
$ws.Range("X15").Value = 75.8584487278179
$ws.Range("B16").Value = 0
$ws.Range("H16").Value = 173.7089986416413
$ws.Range("K16").Value = 22.2359132871029
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 0
$ws.Range("W16").Value = 0
$ws.Range("X16").Value = 0
$ws.Range("AA16").Value = 0
$ws.Range("C17").Value = 0
$ws.Range("I17").Value = 181.050703490946
$ws.Range("M17").Value = 0
$ws.Range("O17").Value = 191.1339795008002
$ws.Range("R17").Value = 0
$ws.Range("S17").Value = 0
$ws.Range("V17").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("J18").Value = 27.67436238789873
$ws.Range("K18").Value = 0
$ws.Range("AA18").Value = 3.477124494654427
$ws.Range("B19").Value = 48.96766975713056
$ws.Range("I19").Value = 135.3966807924818
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("V19").Value = 141.5892427938283
$ws.Range("W19").Value = 98.49183067922618
$ws.Range("Y19").Value = 144.8837424607011
$ws.Range("C20").Value = 83.05115025630701
$ws.Range("G20").Value = 0
$ws.Range("M20").Value = 5.810390817416455
$ws.Range("N20").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("V20").Value = 0
$ws.Range("Z20").Value = 0
$ws.Range("AA20").Value = 164.1831165485813
$ws.Range("AD20").Value = 109.8081415153659
$ws.Range("O21").Value = 113.817046644803
$ws.Range("S21").Value = 0
$ws.Range("X21").Value = 0
$ws.Range("AA21").Value = 0
$ws.Range("J22").Value = 73.98213654092689
$ws.Range("Q22").Value = 0.3050613965543514
$ws.Range("Z22").Value = 159.9123534890697
$ws.Range("AA22").Value = 0
$ws.Range("AE22").Value = 105.4185342872002
$ws.Range("E23").Value = 0
$ws.Range("P23").Value = 181.1219945622685
$ws.Range("X23").Value = 55.14463800256786
$ws.Range("I24").Value = 79.30013583558666
$ws.Range("V24").Value = 188.174734618735
$ws.Range("AC24").Value = 12.22364148619475
$ws.Range("AE24").Value = 0
$ws.Range("B25").Value = 0
$ws.Range("T25").Value = 107.2530760813871
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 102.9108913985278
$ws.Range("V26").Value = 36.79175731258779
$ws.Range("AE26").Value = 41.7294991912611
$ws.Range("B27").Value = 0
$ws.Range("I27").Value = 162.7638769258638
$ws.Range("Q27").Value = 0
$ws.Range("AB27").Value = 93.60126849868251
$ws.Range("B28").Value = 131.6757503849493
$ws.Range("L28").Value = 160.6846614163516
$ws.Range("Z28").Value = 6.152500669787653
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("J29").Value = 45.32619768342816
$ws.Range("R29").Value = 0
$ws.Range("T29").Value = 54.69283743325095
$ws.Range("Z29").Value = 0
$ws.Range("AD29").Value = 151.3117526103657
$ws.Range("P30").Value = 127.7914282145141
$ws.Range("U30").Value = 49.65861353749552
$ws.Range("Y30").Value = 109.8278416991078
$ws.Range("Z30").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("H31").Value = 127.1312804989518
$ws.Range("Q31").Value = 150.7790807799552
$ws.Range("R31").Value = 43.32180468213116
$ws.Range("Y31").Value = 0
$ws.Range("AC31").Value = 42.05273010064246
